$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.233.69'
$ws.Range("E2").Value = '  +1.50%  '

$ws.Range("D3").Value = '2.491.59'
$ws.Range("E3").Value = '  +0.78%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '322.18'
$ws.Range("E5").Value = '  +0.07%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.87'
$ws.Range("E6").Value = '  +2.94%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.525'
$ws.Range("E7").Value = '  +0.82%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.536'
$ws.Range("E9").Value = '  -0.75%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.91'
$ws.Range("E10").Value = '  +7.37%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0813'
$ws.Range("E11").Value = '  -0.18%  '

$ws.Range("E12").Value = '  +0.52%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.34'
$ws.Range("E13").Value = '  -0.13%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.17'
$ws.Range("E14").Value = '  +0.83%  '

$ws.Range("D15").Value = '2.879.73'
$ws.Range("E15").Value = '  +0.49%  '

$ws.Range("D16").Value = '2.492.25'
$ws.Range("E16").Value = '  +0.47%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.848'
$ws.Range("E17").Value = '  +0.08%  '

$ws.Range("D18").Value = '47.140.84'
$ws.Range("E18").Value = '  +1.48%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.77'
$ws.Range("E19").Value = '  +1.00%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.63'
$ws.Range("E20").Value = '  +2.25%  '

$ws.Range("E21").Value = '  -0.08%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.75'
$ws.Range("E22").Value = '  +15.69%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.66'
$ws.Range("E23").Value = '  +0.11%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '246.85'
$ws.Range("E24").Value = '  -0.74%  '

$ws.Range("E25").Value = '  +1.43%  '

$ws.Range("E26").Value = '  -0.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.88'
$ws.Range("E27").Value = '  -1.12%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.28'
$ws.Range("E28").Value = '  +3.69%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.06'
$ws.Range("E29").Value = '  +2.75%  '

$ws.Range("E30").Value = '  +9.36%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.40'
$ws.Range("E31").Value = '  +1.96%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.93'
$ws.Range("E32").Value = '  +0.46%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.09'
$ws.Range("E33").Value = '  +1.73%  '

$ws.Range("E34").Value = '  +1.63%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0786'
$ws.Range("E35").Value = '  +2.30%  '

$ws.Range("E36").Value = '  +0.18%  '

$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.96'
$ws.Range("E37").Value = '  +2.85%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.69'
$ws.Range("E38").Value = '  +2.20%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.97'
$ws.Range("E39").Value = '  +0.65%  '

$ws.Range("E40").Value = '  +0.48%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '120.70'
$ws.Range("E42").Value = '  -2.29%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.17'
$ws.Range("E43").Value = '  +1.97%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0297'
$ws.Range("E44").Value = '  +0.94%  '

$ws.Range("D45").Value = '1.993.14'
$ws.Range("E45").Value = '  +0.41%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.05'
$ws.Range("E46").Value = '  +2.34%  '

$ws.Range("E47").Value = '  -2.81%  '

$ws.Range("E48").Value = '  -0.82%  '

$ws.Range("E49").Value = '  +0.36%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.12'
$ws.Range("E50").Value = '  -1.88%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '56.86'
$ws.Range("E51").Value = '  +3.78%  '
